$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.045921100599999
$ws.Range("C2").Value = 1.6535763914
$ws.Range("D2").Value = 0.4917327084
$ws.Range("E2").Value = 0.0095544033
$ws.Range("F2").Value = 0.1981495423
$ws.Range("G2").Value = 0.0039626853

$ws.Range("B3").Value = 10.3455890994
$ws.Range("C3").Value = 2.3428745331
$ws.Range("D3").Value = 0.4871404119
$ws.Range("E3").Value = 0.0235880626
$ws.Range("F3").Value = 0.1993579093
$ws.Range("G3").Value = 0.0131629743

$ws.Range("B4").Value = 7.9015546244
$ws.Range("C4").Value = 1.596426723
$ws.Range("D4").Value = 0.4940488978
$ws.Range("E4").Value = 0.0115855836
$ws.Range("F4").Value = 0.1938878187
$ws.Range("G4").Value = 0.0069952424
